$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.549.86"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.640.24"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.37"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  +4.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.85"
$ws.Range("E8").Value = "  -4.50%  "
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").Value = "1.872.57"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "1.641.25"
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.561"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "27.521.77"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.90"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").Value = "0.0₃0722"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.59"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.97"
$ws.Range("E23").Value = "  +7.15%  "
$ws.Range("E24").Value = "  -3.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.09"
$ws.Range("E25").Value = "  +1.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.94"
$ws.Range("E26").Value = "  -3.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.112"
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0484"
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  -3.64%  "
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.907"
$ws.Range("E40").Value = "  +15.90%  "
$ws.Range("E41").Value = "  -2.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.47"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.49"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.26"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.98"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "1.782.14"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("E48").Value = "  -2.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.16"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0982"
$ws.Range("E51").Value = "  -2.90%  "
